$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new columns, copying the header style (bold/centered/bordered)
# from an existing header cell (AC1) so the new headers match the look of the rest.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the team record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 98
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Host "done"
